# Refresh Universalis market-price snapshot columns (currentAveragePrice*,
# LevePrice*/LeveProfit*) across all job leve-profit sheets, per scheduled run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1093.6842
$ws.Range("I4").Value = 797.8333
$ws.Range("J4").Value = 1600.8572
$ws.Range("K4").Value = 797.8333
$ws.Range("L4").Value = 1600.8572
$ws.Range("M4").Value = -683.8333
$ws.Range("N4").Value = -1828.8572
$ws.Range("H33").Value = 311.2353
$ws.Range("I33").Value = 303.2143
$ws.Range("J33").Value = 348.66666
$ws.Range("K33").Value = 303.2143
$ws.Range("L33").Value = 348.66666
$ws.Range("M33").Value = -74.21429999999998
$ws.Range("N33").Value = -806.66666
$ws.Range("H43").Value = 2596.3333
$ws.Range("I43").Value = 3400
$ws.Range("J43").Value = 989
$ws.Range("K43").Value = 3400
$ws.Range("L43").Value = 989
$ws.Range("M43").Value = -3331
$ws.Range("N43").Value = -1127
$ws.Range("H64").Value = 3959.6667
$ws.Range("J64").Value = 3939.5
$ws.Range("L64").Value = 3939.5
$ws.Range("N64").Value = -4435.5
$ws.Range("H67").Value = 3959.6667
$ws.Range("J67").Value = 3939.5
$ws.Range("L67").Value = 3939.5
$ws.Range("N67").Value = -5655.5
$ws.Range("H99").Value = 3478.3333
$ws.Range("J99").Value = 3893.25
$ws.Range("L99").Value = 11679.75
$ws.Range("N99").Value = -14675.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2200.5757
$ws.Range("I2").Value = 2447.6086
$ws.Range("J2").Value = 1632.4
$ws.Range("K2").Value = 2447.6086
$ws.Range("L2").Value = 1632.4
$ws.Range("M2").Value = -2334.6086
$ws.Range("N2").Value = -1858.4
$ws.Range("H19").Value = 37495
$ws.Range("I19").Value = 35000
$ws.Range("J19").Value = 39990
$ws.Range("K19").Value = 35000
$ws.Range("L19").Value = 39990
$ws.Range("M19").Value = -34771
$ws.Range("N19").Value = -40448
$ws.Range("H32").Value = 4118.1396
$ws.Range("I32").Value = 3352
$ws.Range("J32").Value = 14333.333
$ws.Range("K32").Value = 3352
$ws.Range("L32").Value = 14333.333
$ws.Range("M32").Value = -3065
$ws.Range("N32").Value = -14907.333
$ws.Range("H116").Value = 2200.5757
$ws.Range("I116").Value = 2447.6086
$ws.Range("J116").Value = 1632.4
$ws.Range("K116").Value = 2447.6086
$ws.Range("L116").Value = 1632.4
$ws.Range("M116").Value = -153.6086
$ws.Range("N116").Value = -6220.4
$ws.Range("H129").Value = 99400
$ws.Range("J129").Value = 99400
$ws.Range("L129").Value = 99400
$ws.Range("N129").Value = -109400
$ws.Range("H132").Value = 3453122.2
$ws.Range("I132").Value = 4851.409
$ws.Range("K132").Value = 14554.227
$ws.Range("M132").Value = -12024.227

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2200.5757
$ws.Range("I3").Value = 2447.6086
$ws.Range("J3").Value = 1632.4
$ws.Range("K3").Value = 2447.6086
$ws.Range("L3").Value = 1632.4
$ws.Range("M3").Value = -2333.6086
$ws.Range("N3").Value = -1860.4
$ws.Range("H80").Value = 1060.5714
$ws.Range("I80").Value = 1384.5
$ws.Range("J80").Value = 931
$ws.Range("K80").Value = 1384.5
$ws.Range("L80").Value = 931
$ws.Range("M80").Value = -386.5
$ws.Range("N80").Value = -2927
$ws.Range("H82").Value = 3992.4
$ws.Range("I82").Value = 3992.4
$ws.Range("K82").Value = 3992.4
$ws.Range("M82").Value = -3609.4
$ws.Range("H83").Value = 1060.5714
$ws.Range("I83").Value = 1384.5
$ws.Range("J83").Value = 931
$ws.Range("K83").Value = 6922.5
$ws.Range("L83").Value = 4655
$ws.Range("M83").Value = -1930.5
$ws.Range("N83").Value = -14639
$ws.Range("H85").Value = 3992.4
$ws.Range("I85").Value = 3992.4
$ws.Range("K85").Value = 3992.4
$ws.Range("M85").Value = -2666.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4666.3335
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 5499.5
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 5499.5
$ws.Range("M3").Value = -2887
$ws.Range("N3").Value = -5725.5
$ws.Range("H7").Value = 331.66666
$ws.Range("I7").Value = 95
$ws.Range("K7").Value = 95
$ws.Range("M7").Value = 18
$ws.Range("H31").Value = 28574208
$ws.Range("I31").Value = 58826550
$ws.Range("J31").Value = 2548.8333
$ws.Range("K31").Value = 58826550
$ws.Range("L31").Value = 2548.8333
$ws.Range("M31").Value = -58826255
$ws.Range("N31").Value = -3138.8333
$ws.Range("H34").Value = 28574208
$ws.Range("I34").Value = 58826550
$ws.Range("J34").Value = 2548.8333
$ws.Range("K34").Value = 58826550
$ws.Range("L34").Value = 2548.8333
$ws.Range("M34").Value = -58826348
$ws.Range("N34").Value = -2952.8333
$ws.Range("H109").Value = 99999
$ws.Range("J109").Value = 99999
$ws.Range("L109").Value = 99999
$ws.Range("N109").Value = -102079
$ws.Range("H122").Value = 12280.1875
$ws.Range("I122").Value = 16754.3
$ws.Range("J122").Value = 4823.3335
$ws.Range("K122").Value = 50262.89999999999
$ws.Range("L122").Value = 14470.0005
$ws.Range("M122").Value = -47812.89999999999
$ws.Range("N122").Value = -19370.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 12465.556
$ws.Range("I14").Value = 12465.556
$ws.Range("K14").Value = 37396.66800000001
$ws.Range("M14").Value = -37223.66800000001
$ws.Range("H92").Value = 113
$ws.Range("J92").Value = 113
$ws.Range("L92").Value = 339
$ws.Range("N92").Value = -2835
$ws.Range("H122").Value = 49578.145
$ws.Range("J122").Value = 2862.25
$ws.Range("L122").Value = 25760.25
$ws.Range("N122").Value = -30660.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1685705.1
$ws.Range("I113").Value = 2292.75
$ws.Range("J113").Value = 6174804.5
$ws.Range("K113").Value = 2292.75
$ws.Range("L113").Value = 6174804.5
$ws.Range("M113").Value = -122.75
$ws.Range("N113").Value = -6179144.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -205
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -393
$ws.Range("H68").Value = 6946790.5
$ws.Range("I68").Value = 13890347
$ws.Range("J68").Value = 3234
$ws.Range("K68").Value = 13890347
$ws.Range("L68").Value = 3234
$ws.Range("M68").Value = -13889598
$ws.Range("N68").Value = -4732
$ws.Range("H71").Value = 6946790.5
$ws.Range("I71").Value = 13890347
$ws.Range("J71").Value = 3234
$ws.Range("K71").Value = 69451735
$ws.Range("L71").Value = 16170
$ws.Range("M71").Value = -69447991
$ws.Range("N71").Value = -23658
$ws.Range("H93").Value = 1987571.8
$ws.Range("I93").Value = 1522.6842
$ws.Range("J93").Value = 6180342
$ws.Range("K93").Value = 1522.6842
$ws.Range("L93").Value = 6180342
$ws.Range("M93").Value = -274.6841999999999
$ws.Range("N93").Value = -6182838

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 8979.666999999999
$ws.Range("I96").Value = 8097
$ws.Range("J96").Value = 10745
$ws.Range("K96").Value = 8097
$ws.Range("L96").Value = 10745
$ws.Range("M96").Value = -6724
$ws.Range("N96").Value = -13491
$ws.Range("H132").Value = 298946.78
$ws.Range("I132").Value = 4190.8096
$ws.Range("J132").Value = 775091.0600000001
$ws.Range("K132").Value = 12572.4288
$ws.Range("L132").Value = 2325273.18
$ws.Range("M132").Value = -10042.4288
$ws.Range("N132").Value = -2330333.18
